$d = $word.ActiveDocument

$wdFindWrapContinue = 1
$wdReplaceAll = 2

function Replace-InParagraph($paraIndex, $find, $replace) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    if ($rng.Text -notlike "*$find*") {
        Write-Output "WARNING: paragraph $paraIndex does not contain '$find' (text='$($rng.Text)')"
    }
    $ok = $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, $wdFindWrapContinue, $false, $replace, $wdReplaceAll)
    if (-not $ok) {
        Write-Output "WARNING: replace failed for '$find' in paragraph $paraIndex"
    }
}

# English -> Russian menu-item translations (Menu-Languages.docx, Russian)
Replace-InParagraph 8  "Hive Teams" "Команды Hive"
Replace-InParagraph 9  "Roadmap" "План развития"
Replace-InParagraph 10 "SmartHive Discussion" "SmartHive Обсуждения"
Replace-InParagraph 11 "SmartHive Voting" "SmartHive Голосование"
Replace-InParagraph 12 "Publications Archive" "Архив публикаций"
Replace-InParagraph 14 "Meetup" "Meetups"
Replace-InParagraph 15 "Businesses" "Бизнес"
Replace-InParagraph 16 "Benefit for Merchants" "Преимущества для Бизнеса"
Replace-InParagraph 18 "The Other Side" "Другая сторона"
Replace-InParagraph 19 "Resources" "Ресурсы"
Replace-InParagraph 20 "Insight Explorer" "Insight обозреватель"
Replace-InParagraph 21 "Pool overview" "Пулы для майнинга"
Replace-InParagraph 23 "Services" "Сервисы"
Replace-InParagraph 24 "Projects" "Проекты"
Replace-InParagraph 25 "Exchanges" "Биржи"
